# Update workbook/sheet metadata for the new "through" date (12-15 -> 12-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-12-17"

# Update the row label for December
$ws.Range("A13").Value = "December (through 12-17)"

# I12 (2022, November) updated
$ws.Range("I12").Value = 119

# December row (row 13) values for years 2015-2022 (columns B-I)
$ws.Range("B13").Value = 21
$ws.Range("C13").Value = 54
$ws.Range("D13").Value = 67
$ws.Range("E13").Value = 39
$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 80
$ws.Range("H13").Value = 129
$ws.Range("I13").Value = 73

# Total row (row 14) values for years 2015-2022 (columns B-I)
$ws.Range("B14").Value = 312
$ws.Range("C14").Value = 617
$ws.Range("D14").Value = 888
$ws.Range("E14").Value = 721
$ws.Range("F14").Value = 561
$ws.Range("G14").Value = 1344
$ws.Range("H14").Value = 1772
$ws.Range("I14").Value = 1590
